$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on the Price column so numeric-looking strings
# (e.g. "0.4810", "315.57") are stored verbatim instead of being
# auto-converted to floating point numbers (which would drop
# trailing zeros / change representation).
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "28.026.00"
$ws.Range("E2").Value = "  +1.94%  "
$ws.Range("D3").Value = "1.912.72"
$ws.Range("E3").Value = "  +2.62%  "
$ws.Range("E4").Value = "  -0.41%  "
$ws.Range("D5").Value = "315.57"
$ws.Range("E5").Value = "  +1.33%  "
$ws.Range("E6").Value = "  -0.44%  "
$ws.Range("D7").Value = "0.4810"
$ws.Range("E7").Value = "  +0.68%  "
$ws.Range("D8").Value = "0.3814"
$ws.Range("E8").Value = "  +0.00%  "
$ws.Range("D9").Value = "0.07356"
$ws.Range("E9").Value = "  +0.28%  "
$ws.Range("D10").Value = "0.9332"
$ws.Range("E10").Value = "  -0.32%  "
$ws.Range("D11").Value = "20.79"
$ws.Range("E11").Value = "  -0.03%  "
$ws.Range("D12").Value = "0.07802"
$ws.Range("E12").Value = "  -0.09%  "
$ws.Range("D13").Value = "1.867.96"
$ws.Range("E13").Value = "  +0.12%  "
$ws.Range("E14").Value = "  +1.01%  "
$ws.Range("D15").Value = "6.628"
$ws.Range("D16").Value = "91.87"
$ws.Range("E16").Value = "  +1.48%  "
$ws.Range("E17").Value = "  -0.44%  "
$ws.Range("D18").Value = "0.000008879"
$ws.Range("E18").Value = "  +1.07%  "
$ws.Range("E19").Value = "  -0.34%  "
$ws.Range("D20").Value = "28.047.19"
$ws.Range("E20").Value = "  +1.87%  "
$ws.Range("D21").Value = "14.75"
$ws.Range("E21").Value = "  +0.52%  "
$ws.Range("D22").Value = "5.170"
$ws.Range("E22").Value = "  +0.94%  "
$ws.Range("D23").Value = "2.146.88"
$ws.Range("E23").Value = "  +2.34%  "
$ws.Range("E24").Value = "  +1.85%  "
$ws.Range("D25").Value = "156.56"
$ws.Range("E25").Value = "  +1.19%  "
$ws.Range("D26").Value = "1.911"
$ws.Range("E26").Value = "  -1.52%  "
$ws.Range("E27").Value = "  +0.04%  "
$ws.Range("D28").Value = "2.128"
$ws.Range("E28").Value = "  +5.38%  "
$ws.Range("D29").Value = "116.70"
$ws.Range("E29").Value = "  +1.18%  "
$ws.Range("D30").Value = "4.963"
$ws.Range("E30").Value = "  +0.48%  "
$ws.Range("D31").Value = "0.08947"
$ws.Range("E31").Value = "  +0.62%  "
$ws.Range("E32").Value = "  -0.65%  "
$ws.Range("D33").Value = "1.254"
$ws.Range("E33").Value = "  +2.96%  "
$ws.Range("D34").Value = "0.7743"
$ws.Range("E34").Value = "  +2.03%  "
$ws.Range("D35").Value = "4.665"
$ws.Range("E35").Value = "  +1.32%  "
$ws.Range("D36").Value = "2.610"
$ws.Range("D37").Value = "0.02048"
$ws.Range("E37").Value = "  -0.32%  "
$ws.Range("E38").Value = "  -1.06%  "
$ws.Range("D39").Value = "0.5523"
$ws.Range("E39").Value = "  -0.93%  "
$ws.Range("D40").Value = "0.05303"
$ws.Range("E40").Value = "  +0.46%  "
$ws.Range("D41").Value = "2.997"
$ws.Range("E41").Value = "  +0.17%  "
$ws.Range("D42").Value = "7.029"
$ws.Range("E42").Value = "  -0.45%  "
$ws.Range("D43").Value = "0.1527"
$ws.Range("E43").Value = "  +0.05%  "
$ws.Range("D44").Value = "8.499"
$ws.Range("E44").Value = "  -1.95%  "
$ws.Range("D45").Value = "10.72"
$ws.Range("E45").Value = "  -0.07%  "
$ws.Range("D46").Value = "108.64"
$ws.Range("E46").Value = "  +5.43%  "
$ws.Range("E47").Value = "  -1.67%  "
$ws.Range("D48").Value = "1.005"
$ws.Range("E48").Value = "  -0.45%  "
$ws.Range("D49").Value = "1.646"
$ws.Range("E49").Value = "  -0.66%  "
$ws.Range("D50").Value = "67.92"
$ws.Range("E50").Value = "  +0.69%  "
$ws.Range("D51").Value = "0.06078"
$ws.Range("E51").Value = "  -0.17%  "

# Restore default ("Normal") style on the Price column so the saved
# cell style index matches the original workbook (only the cell
# text content should differ, not the formatting/style).
$ws.Range("D2:D51").Style = "Normal"

